$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Label Encode" section (rows 14-22) ---

# Row 14: section title
$ws.Range("B14").Value = "Label Encode"

# Row 15: column headers (mirrors row 3's headers, with a couple renamed).
# Written in B, F, A, C, D, E order so new shared-string entries land at the
# same indices as the target workbook (8 = "Train Error", 9 = "Epoch Needed").
$ws.Range("B15").Value = "Train Error"
$ws.Range("F15").Value = "Epoch Needed"
$ws.Range("A15").Value = "Scaling Factor"
$ws.Range("C15").Value = "Test Error"
$ws.Range("D15").Value = "Maximum % error"
$ws.Range("E15").Value = "Average % error"

# Rows 16-22: data (scientific-notation literals rewritten as plain decimals
# to avoid PowerShell parser issues; the resulting doubles are bit-identical)
$ws.Range("A16").Value = 0.5
$ws.Range("B16").Value = 0.000137
$ws.Range("C16").Value = 0.000104
$ws.Range("D16").Value = 34.959778
$ws.Range("E16").Value = 0.567899
$ws.Range("F16").Value = 7900

$ws.Range("A17").Value = 1.5
$ws.Range("B17").Value = 0.00298
$ws.Range("C17").Value = 0.0001002
$ws.Range("D17").Value = 33.224969
$ws.Range("E17").Value = 0.582593
$ws.Range("F17").Value = 22500

$ws.Range("A18").Value = 3
$ws.Range("B18").Value = 0.000442
$ws.Range("C18").Value = 0.0000863
$ws.Range("D18").Value = 14.1225
$ws.Range("E18").Value = 0.586575
$ws.Range("F18").Value = 17900

$ws.Range("A19").Value = 5
$ws.Range("B19").Value = 0.000735
$ws.Range("C19").Value = 0.0001476
$ws.Range("D19").Value = 54.525543
$ws.Range("E19").Value = 0.600287
$ws.Range("F19").Value = 11000

$ws.Range("A20").Value = 8
$ws.Range("B20").Value = 0.000347
$ws.Range("C20").Value = 0.0001073
$ws.Range("D20").Value = 36.772277
$ws.Range("E20").Value = 0.586688
$ws.Range("F20").Value = 70600

$ws.Range("A21").Value = 10
$ws.Range("B21").Value = 0.001551
$ws.Range("C21").Value = 0.0003197
$ws.Range("D21").Value = 106.27934
$ws.Range("E21").Value = 0.636838
$ws.Range("F21").Value = 20700

$ws.Range("A22").Value = 15
$ws.Range("B22").Value = 0.000458
$ws.Range("C22").Value = 0.0004995
$ws.Range("D22").Value = 150.497
$ws.Range("E22").Value = 0.611398
$ws.Range("F22").Value = 65900

# C18 carries a scientific-notation number format in the target workbook
$ws.Range("C18").NumberFormat = "0.00E+00"

# Update selection / scroll position to match the post-edit view state
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F23").Select()
